$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ row = 13; foreign = "abeille";   grammar = "nf"; pronunciation = "abej";    meaning = "včela" },
    @{ row = 14; foreign = "chat";      grammar = "nm"; pronunciation = "ša";      meaning = "kočka, kocour" },
    @{ row = 15; foreign = "chien";     grammar = "nm"; pronunciation = "šje~";    meaning = "pes" },
    @{ row = 16; foreign = "coq";       grammar = "nm"; pronunciation = "kok";     meaning = "kohout" },
    @{ row = 17; foreign = "loup";      grammar = "nm"; pronunciation = "lu";      meaning = "vlk" },
    @{ row = 18; foreign = "mouton";    grammar = "nm"; pronunciation = "muto~";   meaning = "skopec, beránek" },
    @{ row = 19; foreign = "oiseau";    grammar = "nm"; pronunciation = "u^azo."; meaning = "pták" },
    @{ row = 20; foreign = "perroquet"; grammar = "nm"; pronunciation = "peroke"; meaning = "papoušek" },
    @{ row = 21; foreign = "vache";     grammar = "nf"; pronunciation = "vaš";     meaning = "kráva" }
)

foreach ($r in $rows) {
    $rowNum = $r.row
    $ws.Cells.Item($rowNum, 1).Value = $r.foreign
    $ws.Cells.Item($rowNum, 2).Value = $r.grammar
    $ws.Cells.Item($rowNum, 3).Value = $r.pronunciation
    $ws.Cells.Item($rowNum, 4).Value = $r.meaning
}

# Enter the formula once across the whole block so Excel records it as a
# single shared formula (F13:F21), matching how the original F3:F11 block
# was authored.
$f = ' "{ ""foreign"": """ & A13 & """, ""grammar"": """ & B13 & """, ""pronunciation"": """ & C13 & """, ""meaning"": """ & D13 & """ },"'
$ws.Range("F13:F21").Formula = "=" + $f

[void]$ws.Range("F13").Select()
